# Auto-generated edit script applying the Leviathan_Profits market-data refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to match the latest scrape.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3898.4
$ws.Range("I76").Value = 3874.5
$ws.Range("K76").Value = 3874.5
$ws.Range("M76").Value = -3559.5
$ws.Range("H79").Value = 3898.4
$ws.Range("I79").Value = 3874.5
$ws.Range("K79").Value = 3874.5
$ws.Range("M79").Value = -2782.5
$ws.Range("H97").Value = 1813.1666
$ws.Range("I97").Value = 777
$ws.Range("K97").Value = 2331
$ws.Range("M97").Value = -1835
$ws.Range("H98").Value = 2120.9412
$ws.Range("I98").Value = 1075.5714
$ws.Range("K98").Value = 1075.5714
$ws.Range("M98").Value = 422.4286
$ws.Range("H101").Value = 1142.6666
$ws.Range("I101").Value = 1203.4166
$ws.Range("K101").Value = 3610.2498
$ws.Range("M101").Value = -1988.2498
$ws.Range("H122").Value = 2120.9412
$ws.Range("I122").Value = 1075.5714
$ws.Range("K122").Value = 3226.7142
$ws.Range("M122").Value = -776.7142000000003
$ws.Range("H132").Value = 2986.6829
$ws.Range("I132").Value = 1742.5555
$ws.Range("K132").Value = 5227.666499999999
$ws.Range("M132").Value = -2697.666499999999
$ws.Range("H137").Value = 1184
$ws.Range("I137").Value = 1121.5
$ws.Range("J137").Value = 1309
$ws.Range("K137").Value = 3364.5
$ws.Range("L137").Value = 3927
$ws.Range("M137").Value = -814.5
$ws.Range("N137").Value = -9027
$ws.Range("H138").Value = 2327.919
$ws.Range("J138").Value = 2612.3809
$ws.Range("L138").Value = 7837.1427
$ws.Range("N138").Value = -18117.1427
$ws.Range("H141").Value = 2343
$ws.Range("I141").Value = 2343
$ws.Range("K141").Value = 7029
$ws.Range("M141").Value = -1849

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14925.913
$ws.Range("I32").Value = 4242.316
$ws.Range("J32").Value = 65673
$ws.Range("K32").Value = 4242.316
$ws.Range("L32").Value = 65673
$ws.Range("M32").Value = -3955.316
$ws.Range("N32").Value = -66247
$ws.Range("H61").Value = 1914.1666
$ws.Range("I61").Value = 1802.7693
$ws.Range("K61").Value = 1802.7693
$ws.Range("M61").Value = -1590.7693
$ws.Range("H74").Value = 1874.25
$ws.Range("I74").Value = 1832.5333
$ws.Range("K74").Value = 1832.5333
$ws.Range("M74").Value = -958.5333000000001
$ws.Range("H77").Value = 1874.25
$ws.Range("I77").Value = 1832.5333
$ws.Range("K77").Value = 9162.666499999999
$ws.Range("M77").Value = -4794.666499999999
$ws.Range("H80").Value = 15666.667
$ws.Range("I80").Value = 7000
$ws.Range("K80").Value = 7000
$ws.Range("M80").Value = -6002
$ws.Range("H83").Value = 15666.667
$ws.Range("I83").Value = 7000
$ws.Range("K83").Value = 21000
$ws.Range("M83").Value = -16008
$ws.Range("H92").Value = 29516.334
$ws.Range("J92").Value = 29516.334
$ws.Range("L92").Value = 29516.334
$ws.Range("N92").Value = -34508.334
$ws.Range("H97").Value = 2095.739
$ws.Range("I97").Value = 1941.1364
$ws.Range("J97").Value = 5497
$ws.Range("K97").Value = 1941.1364
$ws.Range("L97").Value = 5497
$ws.Range("M97").Value = -1445.1364
$ws.Range("N97").Value = -6489
$ws.Range("H136").Value = 1914.1666
$ws.Range("I136").Value = 1802.7693
$ws.Range("K136").Value = 5408.3079
$ws.Range("M136").Value = -2858.3079

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2486.7693
$ws.Range("I86").Value = 2569.7778
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 2569.7778
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -1446.7778
$ws.Range("N86").Value = -4546
$ws.Range("H89").Value = 2486.7693
$ws.Range("I89").Value = 2569.7778
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 12848.889
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -7232.888999999999
$ws.Range("N89").Value = -22732
$ws.Range("H132").Value = 133783.5
$ws.Range("J132").Value = 133783.5
$ws.Range("L132").Value = 133783.5
$ws.Range("N132").Value = -143903.5
$ws.Range("H140").Value = 87500
$ws.Range("J140").Value = 87500
$ws.Range("L140").Value = 87500
$ws.Range("N140").Value = -97860

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10305.596
$ws.Range("I31").Value = 2994.7942
$ws.Range("K31").Value = 2994.7942
$ws.Range("M31").Value = -2699.7942
$ws.Range("H34").Value = 10305.596
$ws.Range("I34").Value = 2994.7942
$ws.Range("K34").Value = 2994.7942
$ws.Range("M34").Value = -2792.7942
$ws.Range("H99").Value = 29779.8
$ws.Range("I99").Value = 54859.6
$ws.Range("J99").Value = 4700
$ws.Range("K99").Value = 54859.6
$ws.Range("L99").Value = 4700
$ws.Range("M99").Value = -53361.6
$ws.Range("N99").Value = -7696
$ws.Range("H126").Value = 29779.8
$ws.Range("I126").Value = 54859.6
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 164578.8
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -162108.8
$ws.Range("N126").Value = -19040

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1612.7059
$ws.Range("I14").Value = 1612.7059
$ws.Range("K14").Value = 4838.1177
$ws.Range("M14").Value = -4665.1177
$ws.Range("H98").Value = 615.1429000000001
$ws.Range("I98").Value = 816.5
$ws.Range("K98").Value = 2449.5
$ws.Range("M98").Value = -951.5
$ws.Range("H107").Value = 633.8214
$ws.Range("I107").Value = 670
$ws.Range("J107").Value = 627.7917
$ws.Range("K107").Value = 2010
$ws.Range("L107").Value = 1883.3751
$ws.Range("M107").Value = -90
$ws.Range("N107").Value = -5723.3751
$ws.Range("H110").Value = 13013
$ws.Range("I110").Value = 9026
$ws.Range("K110").Value = 27078
$ws.Range("M110").Value = -22988
$ws.Range("H113").Value = 2171
$ws.Range("J113").Value = 2699.8
$ws.Range("L113").Value = 8099.400000000001
$ws.Range("N113").Value = -12439.4
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1775.9048
$ws.Range("I132").Value = 1042.5714
$ws.Range("J132").Value = 2142.5715
$ws.Range("K132").Value = 9383.142600000001
$ws.Range("L132").Value = 19283.1435
$ws.Range("M132").Value = -6853.142600000001
$ws.Range("N132").Value = -24343.1435
$ws.Range("H133").Value = 4324.75
$ws.Range("I133").Value = 1319.6
$ws.Range("J133").Value = 9333.333000000001
$ws.Range("K133").Value = 3958.8
$ws.Range("L133").Value = 27999.999
$ws.Range("M133").Value = 1101.2
$ws.Range("N133").Value = -38119.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2792.3333
$ws.Range("I80").Value = 1692.875
$ws.Range("J80").Value = 4991.25
$ws.Range("K80").Value = 1692.875
$ws.Range("L80").Value = 4991.25
$ws.Range("M80").Value = -694.875
$ws.Range("N80").Value = -6987.25
$ws.Range("H83").Value = 2792.3333
$ws.Range("I83").Value = 1692.875
$ws.Range("J83").Value = 4991.25
$ws.Range("K83").Value = 8464.375
$ws.Range("L83").Value = 24956.25
$ws.Range("M83").Value = -3472.375
$ws.Range("N83").Value = -34940.25
$ws.Range("H126").Value = 1658.6666
$ws.Range("I126").Value = 999
$ws.Range("K126").Value = 2997
$ws.Range("M126").Value = -527
$ws.Range("H136").Value = 25452.375
$ws.Range("J136").Value = 25452.375
$ws.Range("L136").Value = 76357.125
$ws.Range("N136").Value = -81457.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 72551.14
$ws.Range("J16").Value = 250799.75
$ws.Range("L16").Value = 250799.75
$ws.Range("N16").Value = -251139.75
$ws.Range("H82").Value = 1538.2222
$ws.Range("I82").Value = 1499.5
$ws.Range("J82").Value = 1569.2
$ws.Range("K82").Value = 1499.5
$ws.Range("L82").Value = 1569.2
$ws.Range("M82").Value = -1138.5
$ws.Range("N82").Value = -2291.2
$ws.Range("H85").Value = 1538.2222
$ws.Range("I85").Value = 1499.5
$ws.Range("J85").Value = 1569.2
$ws.Range("K85").Value = 1499.5
$ws.Range("L85").Value = 1569.2
$ws.Range("M85").Value = -251.5
$ws.Range("N85").Value = -4065.2
$ws.Range("H93").Value = 9763.575000000001
$ws.Range("J93").Value = 67391.8
$ws.Range("L93").Value = 67391.8
$ws.Range("N93").Value = -69887.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2177.6667
$ws.Range("I81").Value = 2166.4167
$ws.Range("K81").Value = 4332.8334
$ws.Range("M81").Value = -3271.8334
$ws.Range("H84").Value = 2177.6667
$ws.Range("I84").Value = 2166.4167
$ws.Range("K84").Value = 21664.167
$ws.Range("M84").Value = -16360.167
$ws.Range("H137").Value = 91500
$ws.Range("J137").Value = 91500
$ws.Range("L137").Value = 91500

